# Append: 2025-10-13 12:48 JST
# Update the "取得日時" (acquisition timestamp) column on the ランサーズ sheet
# for all existing data rows from the previous run timestamp to the new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-10-13 12:38:17"
$newTimestamp = "2025-10-13 12:48:29"

# Find last used row in column A and update every matching cell.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value2 = $newTimestamp
    }
}

$wb.Save()
